# Applies the cryptos list refresh (prices / volume% / two coin-row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.863.66'
$ws.Range("E2").Value = '  +0.23%  '

# Row 3
$ws.Range("D3").Value = '2.532.11'
$ws.Range("E3").Value = '  +0.08%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.98'
$ws.Range("E5").Value = '  +0.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.73'
$ws.Range("E6").Value = '  -1.00%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("E8").Value = '  -1.21%  '

# Row 9
$ws.Range("D9").Value = '2.531.64'
$ws.Range("E9").Value = '  +0.11%  '

# Row 10
$ws.Range("E10").Value = '  -1.67%  '

# Row 11
$ws.Range("E11").Value = '  +1.82%  '

# Row 12
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.343'
$ws.Range("E12").Value = '  -0.32%  '

# Row 13
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.03'
$ws.Range("E13").Value = '  -2.83%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.42'
$ws.Range("E14").Value = '  -1.24%  '

# Row 15
$ws.Range("D15").Value = '3.018.02'
$ws.Range("E15").Value = '  +0.92%  '

# Row 16
$ws.Range("B16").Value = 'Binance-PegBSC-USD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.43'
$ws.Range("E16").Value = '  +144.50%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000177'
$ws.Range("E17").Value = '  -0.80%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '67.861.83'
$ws.Range("E18").Value = '  +0.56%  '

# Row 19
$ws.Range("D19").Value = '2.532.94'
$ws.Range("E19").Value = '  +0.20%  '

# Row 20
$ws.Range("E20").Value = '  +2.89%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.01'
$ws.Range("E21").Value = '  -2.15%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '369.73'
$ws.Range("E22").Value = '  +3.14%  '

# Row 23
$ws.Range("E23").Value = '  -1.37%  '

# Row 24
$ws.Range("E24").Value = '  -1.81%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.88'
$ws.Range("E25").Value = '  +2.67%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.26%  '

# Row 27
$ws.Range("E27").Value = '  -4.25%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.95'
$ws.Range("E28").Value = '  -3.02%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0967'
$ws.Range("E30").Value = '  -2.31%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '538.65'
$ws.Range("E31").Value = '  -2.64%  '

# Row 32
$ws.Range("E32").Value = '  +0.23%  '

# Row 33
$ws.Range("E33").Value = '  -2.58%  '

# Row 34
$ws.Range("E34").Value = '  +0.20%  '

# Row 35
$ws.Range("E35").Value = '  -1.48%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.06%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.84'
$ws.Range("E37").Value = '  +1.11%  '

# Row 38
$ws.Range("E38").Value = '  -2.67%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.19'
$ws.Range("E39").Value = '  +2.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.60'
$ws.Range("E40").Value = '  +0.66%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.14'
$ws.Range("E41").Value = '  -0.46%  '

# Row 42
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.350'
$ws.Range("E42").Value = '  -1.65%  '

# Row 43
$ws.Range("E43").Value = '  -2.04%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("E44").Value = '  -2.18%  '

# Row 45
$ws.Range("E45").Value = '  +0.03%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.36'
$ws.Range("E46").Value = '  -1.24%  '

# Row 47
$ws.Range("E47").Value = '  +3.27%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '147.88'
$ws.Range("E48").Value = '  -1.06%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.71'
$ws.Range("E49").Value = '  -0.02%  '

# Row 50
$ws.Range("E50").Value = '  -1.98%  '

# Row 51
$ws.Range("E51").Value = '  +0.75%  '
